# Update Quarto workflow: install R for CI
# (Values below reflect the regenerated RLIe_GET_L3 bootstrap confidence
# intervals that accompanied the workflow change.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F18" = 0.5848484848484848
    "G18" = 0.7075757575757575
    "F19" = 0.5696969696969697
    "G19" = 0.693939393939394
    "F20" = 0.5666666666666667
    "G20" = 0.6893939393939394
    "G21" = 0.693939393939394
    "F26" = 0.7850746268656716
    "F27" = 0.7462686567164178
    "F28" = 0.7343283582089553
    "G29" = 0.8716417910447761
    "G30" = 0.9066666666666666
    "G31" = 0.8853333333333333
    "G32" = 0.8826666666666667
    "F33" = 0.7466666666666666
    "G33" = 0.8746666666666667
    "F38" = 0.8634920634920635
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
